$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A183").Value = "United States"
$ws.Range("A182").Value = "Tanzania"
$ws.Range("A139").Value = "Moldova"
$ws.Range("A140").Value = "Macedonia"
$ws.Range("A138").Value = "South Korea"
